$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.191.47"
$ws.Range("E2").Value = "  -3.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.687.07"
$ws.Range("E3").Value = "  -4.55%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.75"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.13"
$ws.Range("E6").Value = "  +9.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.678.45"
$ws.Range("E7").Value = "  -4.60%  "

$ws.Range("E8").Value = "  -6.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("E10").Value = "  -5.08%  "

$ws.Range("E11").Value = "  -7.37%  "

$ws.Range("E12").Value = "  +4.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000290"
$ws.Range("E13").Value = "  -9.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.35"
$ws.Range("E14").Value = "  -9.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.272.13"
$ws.Range("E15").Value = "  -4.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.684.09"
$ws.Range("E16").Value = "  -5.19%  "

$ws.Range("E17").Value = "  -9.84%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.77"
$ws.Range("E19").Value = "  -7.67%  "

$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.12"
$ws.Range("E20").Value = "  -7.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.959.75"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "408.10"
$ws.Range("E22").Value = "  -6.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.55"
$ws.Range("E23").Value = "  -3.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.64"
$ws.Range("E24").Value = "  -6.01%  "

$ws.Range("E25").Value = "  -7.98%  "

$ws.Range("E26").Value = "  -8.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.93"
$ws.Range("E27").Value = "  -3.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.85"
$ws.Range("E28").Value = "  -5.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.04"
$ws.Range("E29").Value = "  +1.70%  "

$ws.Range("E30").Value = "  -9.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.72"
$ws.Range("E31").Value = "  -6.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.24"
$ws.Range("E32").Value = "  -10.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.43"
$ws.Range("E33").Value = "  -8.27%  "

$ws.Range("E34").Value = "  -6.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "43.42"
$ws.Range("E35").Value = "  -9.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.31"
$ws.Range("E36").Value = "  -8.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "600.94"
$ws.Range("E37").Value = "  -5.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0882"
$ws.Range("E38").Value = "  -9.90%  "

$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("E40").Value = "  -6.40%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("E42").Value = "  -7.44%  "

$ws.Range("E43").Value = "  +2.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.00"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0436"
$ws.Range("E45").Value = "  -7.30%  "

$ws.Range("E46").Value = "  -12.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.18"
$ws.Range("E47").Value = "  -8.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.73"
$ws.Range("E48").Value = "  -3.86%  "

$ws.Range("E49").Value = "  -7.02%  "

$ws.Range("E50").Value = "  -4.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.731.52"
$ws.Range("E51").Value = "  -3.77%  "
